# Update ticket-sale counts (column F) on the "展览" and "全部类型" sheets
# to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1000
$ws1.Range("F6").Value = 2358
$ws1.Range("F7").Value = 201

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1000
$ws4.Range("F8").Value = 2358
$ws4.Range("F10").Value = 201
